# Add the 2022-Q4 sheet (commit: "feat: add 2022-Q4 data")
#
# Before: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4
# After:  总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4
#
# 1) Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q4 and push the existing rows down by one.
# 2) Insert a brand-new worksheet named "2022-Q4" right after "总计"
#    (i.e. before the current "2022-Q3" sheet) and populate it with the
#    per-fund holdings data, matching the layout/styling used by the
#    other quarterly sheets.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# --- 1) shift 总计 rows down one slot (bottom-up, with literal target values
#        so no reads-through-COM / float round-trips are needed) -------------

# row6 <- old row5 (2021-Q4 / 5 / 3.02); A6 needs the bold/border "index
# column" style that A2:A5 already have, copy it from A5 before writing.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 5
$summary.Range("D6").Value = 3.02

# row5 <- old row4 (2022-Q1 / 9 / 6.86)
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 9
$summary.Range("D5").Value = 6.86

# row4 <- old row3 (2022-Q2 / 42 / 10.18)
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 42
$summary.Range("D4").Value = 10.18

# row3 <- old row2 (2022-Q3 / 8 / 0.08)
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 8
$summary.Range("D3").Value = 0.08

# row2 <- brand-new (2022-Q4 / 6 / 0.08)
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.08

# --- 2) insert new "2022-Q4" worksheet before "2022-Q3" ---------------------

$q3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Copy the header/index-column styling (bold + border + center) from the
# summary sheet's header row, which already uses the same style.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A7").PasteSpecial(-4122)

# Headers
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold text in every other quarterly sheet (fund
# codes with leading zeros, and figures that must keep trailing zeros) --
# format them as text before writing so they aren't coerced to numbers.
# (NB: multi-area "B2:B7,D2:G7" ranges only apply a property write to the
# first area in this host, so each column is formatted separately.)
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:D7").NumberFormat = "@"
$q4.Range("E2:E7").NumberFormat = "@"
$q4.Range("F2:F7").NumberFormat = "@"
$q4.Range("G2:G7").NumberFormat = "@"

# Row 2: 588160
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "588160"
$q4.Range("C2").Value = "南方上证科创板新材料ETF"
$q4.Range("D2").Value = "0.90"
$q4.Range("E2").Value = "98.46"
$q4.Range("F2").Value = "3.53"
$q4.Range("G2").Value = "0.0318"
$q4.Range("H2").Value = 8

# Row 3: 588010
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "588010"
$q4.Range("C3").Value = "博时上证科创板新材料ETF"
$q4.Range("D3").Value = "0.53"
$q4.Range("E3").Value = "98.90"
$q4.Range("F3").Value = "3.55"
$q4.Range("G3").Value = "0.0188"
$q4.Range("H3").Value = 8

# Row 4: 015148
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "015148"
$q4.Range("C4").Value = "华安中证1000指数增强A"
$q4.Range("D4").Value = "1.42"
$q4.Range("E4").Value = "92.79"
$q4.Range("F4").Value = "0.93"
$q4.Range("G4").Value = "0.0132"
$q4.Range("H4").Value = 3

# Row 5: 015149
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "015149"
$q4.Range("C5").Value = "华安中证1000指数增强C"
$q4.Range("D5").Value = "1.08"
$q4.Range("E5").Value = "92.79"
$q4.Range("F5").Value = "0.93"
$q4.Range("G5").Value = "0.0100"
$q4.Range("H5").Value = 3

# Row 6: 005000
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "005000"
$q4.Range("C6").Value = "泰康泉林量化价值精选混合A"
$q4.Range("D6").Value = "0.31"
$q4.Range("E6").Value = "89.21"
$q4.Range("F6").Value = "1.75"
$q4.Range("G6").Value = "0.0054"
$q4.Range("H6").Value = 6

# Row 7: 005111
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "005111"
$q4.Range("C7").Value = "泰康泉林量化价值精选混合C"
$q4.Range("D7").Value = "0.14"
$q4.Range("E7").Value = "89.21"
$q4.Range("F7").Value = "1.75"
$q4.Range("G7").Value = "0.0024"
$q4.Range("H7").Value = 6

$summary.Activate()
$summary.Range("A1").Select()
